$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.385.40"
$ws.Range("D3").Value = "1.848.45"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "240.25"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07627"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2906"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.67"
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07732"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.022"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "0.6784"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001058"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.17"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.153"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "29.421.54"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "226.66"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.486"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "157.98"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1378"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.382"
$ws.Range("E27").Value = "  +5.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.463"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05593"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.122"
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.6949"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01799"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "1.228.61"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.720"
$ws.Range("E38").Value = "  -1.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.385"
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9041"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.60"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.89"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.162"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4010"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.982"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "1.679"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "0.1142"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05702"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4628"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.344"
$ws.Range("E51").Value = "  +0.25%  "
